$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1, copying the formatting (bold, border,
# centered) used by the other header cells (e.g. G1 "sum")
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data value for row 2
$ws.Range("H2").Value = 0
